# Update capital structure database values for Cyprus Paper/Forest Products
# rows 2 and 3 (data rows) - same new values applied to both rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D"  = 0.0375
    "G"  = 0.03805668016194332
    "H"  = 0.03805668016194332
    "I"  = -0.03171390013495277
    "J"  = -0.03171390013495277
    "K"  = -0.287
    "L"  = -0.03873144399460188
    "U"  = 0.032
    "V"  = 0.01006289308176101
    "W"  = -0.02094890510948905
    "X"  = 0.1201912746571094
    "Y"  = -0.1411401797665984
    "Z"  = 0.4176765684008794
    "AA" = -0.01324615297897525
    "AB" = 0.08527736457969526
    "AC" = -0.09852351755867052
    "AD" = 3.1
    "AF" = 3.1
    "AG" = 3.068
    "AH" = 0.4936305732484076
    "AI" = 0.1878787878787879
    "AJ" = 0.4910371318822023
    "AK" = 0.1863007043964051
    "AL" = 0.073
    "AM" = 0.064
    "AN" = 17.31843575418995
    "AO" = -3.219178082191781
    "AP" = 17.13966480446927
    "AQ" = -3.671875
}

foreach ($row in 2..3) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
